$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,4).Value = "'247.65"

# Row 3
$ws.Cells.Item(3,4).Value = "'21.79"

# Row 4
$ws.Cells.Item(4,4).Value = "'5.463"

# Row 5
$ws.Cells.Item(5,4).Value = "'0.05695"

# Row 6
$ws.Cells.Item(6,4).Value = "'3.380"

# Row 7
$ws.Cells.Item(7,4).Value = "'0.8037"

# Row 8
$ws.Cells.Item(8,4).Value = "'1.038"

# Row 9
$ws.Cells.Item(9,4).Value = "'0.1474"

# Row 10
$ws.Cells.Item(10,4).Value = "'0.07314"

# Row 11
$ws.Cells.Item(11,4).Value = "'0.03167"

# Row 12
$ws.Cells.Item(12,4).Value = "'0.02935"

# Row 13
$ws.Cells.Item(13,4).Value = "'0.09281"

# Row 14
$ws.Cells.Item(14,4).Value = "'0.001664"

# Row 15
$ws.Cells.Item(15,4).Value = "'3.211"

# Row 16
$ws.Cells.Item(16,4).Value = "'0.04701"

# Row 17
$ws.Cells.Item(17,4).Value = "'0.0005861"

# Row 18
$ws.Cells.Item(18,4).Value = "'0.006349"

# Row 19
$ws.Cells.Item(19,4).Value = "'0.005044"

# Row 22
$ws.Cells.Item(22,4).Value = "'0.0003201"

# Row 23
$ws.Cells.Item(23,4).Value = "'3.772"

# Row 24
$ws.Cells.Item(24,4).Value = "'6.426"

# Row 26
$ws.Cells.Item(26,4).Value = "'0.3289"

# Row 27
$ws.Cells.Item(27,4).Value = "'0.1299"

# Row 40
$ws.Cells.Item(40,4).Value = "'0.04103"

# Row 41
$ws.Cells.Item(41,2).Value = "BKEXToken"
$ws.Cells.Item(41,3).Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Cells.Item(41,4).Value = "'0.1044"
$ws.Cells.Item(41,5).Value = "40BKEXTokenBKK"

# Row 42
$ws.Cells.Item(42,2).Value = "CEJI"
$ws.Cells.Item(42,3).Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Cells.Item(42,4).Value = "'0.002969"
$ws.Cells.Item(42,5).Value = "41CEJICEJI"

# Row 43
$ws.Cells.Item(43,2).Value = "KickToken"
$ws.Cells.Item(43,3).Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Cells.Item(43,4).Value = "'0.006948"
$ws.Cells.Item(43,5).Value = "42KickTokenKICK"

# Row 44
$ws.Cells.Item(44,4).Value = "'0.008037"

# Row 45
$ws.Cells.Item(45,4).Value = "'0.00005830"

# Row 47
$ws.Cells.Item(47,2).Value = "ACDXExchange"
$ws.Cells.Item(47,3).Value = "https://coinranking.com/coin/-y35lbZ7U+acdxexchange-acxt"
$ws.Cells.Item(47,4).Value = "'0.0005801"
$ws.Cells.Item(47,5).Value = "46ACDXExchangeACXTWorstin24h"

# Row 48
$ws.Cells.Item(48,2).Value = "CoinbaseStockToken"
$ws.Cells.Item(48,3).Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Cells.Item(48,4).Value = "'0.6826"
$ws.Cells.Item(48,5).Value = "47CoinbaseStockTokenCOIN"

# Row 49
$ws.Cells.Item(49,2).Value = "BOLO"
$ws.Cells.Item(49,3).Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Cells.Item(49,4).Value = "'0.009752"
$ws.Cells.Item(49,5).Value = "48BOLOBOLO"

# Row 50
$ws.Cells.Item(50,2).Value = "CryptobidCoin"
$ws.Cells.Item(50,3).Value = "https://coinranking.com/coin/h39bvStAP+cryptobidcoin-cbc"
$ws.Cells.Item(50,4).Value = "'0.00002101"
$ws.Cells.Item(50,5).Value = "49CryptobidCoinCBC"

# Row 51
$ws.Cells.Item(51,2).Value = "SpecialPowerGold"
$ws.Cells.Item(51,3).Value = "https://coinranking.com/coin/jPTWzmsWb+specialpowergold-spg"
$ws.Cells.Item(51,4).Value = "'0.01010"
$ws.Cells.Item(51,5).Value = "50SpecialPowerGoldSPG"
